$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows for the Ntf3-Ntrk2 table (row 2 through row 13).
# Each row array: SendingCluster, LigandSymbol, ReceptorSymbol, TargetCluster, then 16 numeric columns (E..T)
$rows = @(
  @("ECs", "Ntf3", "Ntrk2", "ECs", 3, 1, 6.174425666666667, 18.523277, 0.4796039576068858, 0.5423873011998577, 3, 1, 1.098765, 3.296295, 0.04219121913039277, 0.04777771124863077, 6.784242817635, 61.05818535871499, 0.02023507567119573, 0.02591402386165093),
  @("ECs", "Ntf3", "Ntrk2", "FAPs", 3, 1, 6.174425666666667, 18.523277, 0.4796039576068858, 0.5423873011998577, 3, 1, 15.79937666666667, 47.39813, 0.6066765532820465, 0.6870059169052114, 97.55207680800112, 877.96869127201, 0.2909644759413743, 0.3726232851785513),
  @("ECs", "Ntf3", "Ntrk2", "M1", 3, 1, 6.174425666666667, 18.523277, 0.4796039576068858, 0.5423873011998577, 1, 0.3333333333333333, 0.009167, 0.027501, 0.000352001479632415, 0.0003986096017039115, 0.05660096008633334, 0.509408640777, 0.0001688213027151859, 0.0002162007861005347),
  @("ECs", "Ntf3", "Ntrk2", "sCs", 3, 1, 6.174425666666667, 18.523277, 0.4796039576068858, 0.5423873011998577, 2, 1, 9.1351955, 18.270391, 0.3507802261079284, 0.264817762244454, 56.40458556521784, 338.427513391307, 0.1682355846916007, 0.1436337913735549),
  @("FAPs", "Ntf3", "Ntrk2", "ECs", 3, 1, 2.228940333333334, 6.686821, 0.1731349056329954, 0.1957994147469982, 3, 1, 1.098765, 3.296295, 0.04219121913039277, 0.04777771124863077, 2.449081625355, 22.041734628195, 0.007304772742681581, 0.00935484790043298),
  @("FAPs", "Ntf3", "Ntrk2", "FAPs", 3, 1, 2.228940333333334, 6.686821, 0.1731349056329954, 0.1957994147469982, 3, 1, 15.79937666666667, 47.39813, 0.6066765532820465, 0.6870059169052114, 35.21586789385889, 316.94281104473, 0.105036887802238, 0.1345153564577653),
  @("FAPs", "Ntf3", "Ntrk2", "M1", 3, 1, 2.228940333333334, 6.686821, 0.1731349056329954, 0.1957994147469982, 1, 0.3333333333333333, 0.009167, 0.027501, 0.000352001479632415, 0.0003986096017039115, 0.02043269603566667, 0.183894264321, 0.00006094374295883292, 0.00007804752672615994),
  @("FAPs", "Ntf3", "Ntrk2", "sCs", 3, 1, 2.228940333333334, 6.686821, 0.1731349056329954, 0.1957994147469982, 2, 1, 9.1351955, 18.270391, 0.3507802261079284, 0.264817762244454, 20.36180570283517, 122.170834217011, 0.06073230134511696, 0.05185116286207381),
  @("sCs", "Ntf3", "Ntrk2", "ECs", 2, 1, 4.470643, 8.941286, 0.3472611367601188, 0.2618132840531441, 3, 1, 1.098765, 3.296295, 0.04219121913039277, 0.04777771124863077, 4.912186055895, 29.47311633537, 0.01465137071651546, 0.01250883948654686),
  @("sCs", "Ntf3", "Ntrk2", "FAPs", 2, 1, 4.470643, 8.941286, 0.3472611367601188, 0.2618132840531441, 3, 1, 15.79937666666667, 47.39813, 0.6066765532820465, 0.6870059169052114, 70.63337269919667, 423.80023619518, 0.2106751895384342, 0.1798672752688948),
  @("sCs", "Ntf3", "Ntrk2", "M1", 2, 1, 4.470643, 8.941286, 0.3472611367601188, 0.2618132840531441, 1, 0.3333333333333333, 0.009167, 0.027501, 0.000352001479632415, 0.0003986096017039115, 0.040982384381, 0.245894306286, 0.0001222364339583962, 0.0001043612888772168),
  @("sCs", "Ntf3", "Ntrk2", "sCs", 2, 1, 4.470643, 8.941286, 0.3472611367601188, 0.2618132840531441, 2, 1, 9.1351955, 18.270391, 0.3507802261079284, 0.264817762244454, 40.8401978157065, 163.360791262826, 0.1218123400712107, 0.0693328080088252)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowData[$c]
    }
}
